# Rename the sole worksheet from "Property1" to "DataNode" to unify the
# conception of DataNode/DataTable/Entity (per commit message), and restore
# the last-active selected cell on the frozen pane to O40.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"

$ws.Activate()
$ws.Range("O40").Select()
